$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("P1").Value = "betrifft  (0=KL; 1=LP; 2=SuS; 3=alle)"
$ws.Range("Q1").Value = "Unterricht betroffen (0=nein; 1= teilweise; 2=ja)"

$ws.Range("P1:Q1").Font.Size = 12
$ws.Range("P1:Q1").Font.Name = "Calibri (Textkörper)"
$ws.Range("P1:Q1").VerticalAlignment = -4160
$ws.Range("P1:Q1").WrapText = $true
$ws.Range("P1:Q1").Borders.LineStyle = 1
$ws.Range("P1:Q1").Borders.Color = 0

$ws.Range("P2").Value = 2
$ws.Range("Q2").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0
$ws.Range("P4").Value = 3
$ws.Range("Q4").Value = 2
$ws.Range("P5").Value = 0

$ws.Rows.Item(1).RowHeight = 93
$ws.Columns.Item("P").ColumnWidth = 8.26953125

$ws.Range("Q6").Select()
